$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove hyperlink from D13 and clear its value
$ws.Range("D13").Hyperlinks.Delete()
$ws.Range("D13").Value = $null

# Update sheet view: scroll to A7, select D13
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D13").Select()
